# Add a new "time_taken" metadata column (F) to the panel sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: label + same header style (bold/border/centered) as the
# other header cells (copy format from E1, the neighbouring header).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows F2:F8: per-row timestamp strings (kept as plain text, matching
# the other text columns in the sheet).
$ws.Range("F2").Value = "2021-10-05 10:50:26.948622"
$ws.Range("F3").Value = "2021-10-05 10:50:26.948639"
$ws.Range("F4").Value = "2021-10-05 10:50:26.948643"
$ws.Range("F5").Value = "2021-10-05 10:50:26.948647"
$ws.Range("F6").Value = "2021-10-05 10:50:26.948650"
$ws.Range("F7").Value = "2021-10-05 10:50:26.948653"
$ws.Range("F8").Value = "2021-10-05 10:50:26.948656"
